$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1005.25
$ws.Range("I9").Value = 1016.6667
$ws.Range("K9").Value = 1016.6667
$ws.Range("M9").Value = -847.6667
$ws.Range("H88").Value = 1213.9048
$ws.Range("I88").Value = 1422
$ws.Range("J88").Value = 1164.9412
$ws.Range("K88").Value = 1422
$ws.Range("L88").Value = 1164.9412
$ws.Range("M88").Value = -1016
$ws.Range("N88").Value = -1976.9412
$ws.Range("H91").Value = 1213.9048
$ws.Range("I91").Value = 1422
$ws.Range("J91").Value = 1164.9412
$ws.Range("K91").Value = 1422
$ws.Range("L91").Value = 1164.9412
$ws.Range("M91").Value = -18
$ws.Range("N91").Value = -3972.9412
$ws.Range("H92").Value = 54624.2
$ws.Range("I92").Value = 213.07692
$ws.Range("K92").Value = 213.07692
$ws.Range("M92").Value = 1034.92308
$ws.Range("H98").Value = 2549.3809
$ws.Range("J98").Value = 2466.6667
$ws.Range("L98").Value = 2466.6667
$ws.Range("N98").Value = -5462.6667
$ws.Range("H107").Value = 818.72
$ws.Range("I107").Value = 847.4545000000001
$ws.Range("K107").Value = 847.4545000000001
$ws.Range("M107").Value = 1072.5455
$ws.Range("H122").Value = 2549.3809
$ws.Range("J122").Value = 2466.6667
$ws.Range("L122").Value = 7400.000100000001
$ws.Range("N122").Value = -12300.0001
$ws.Range("H135").Value = 1062.6818
$ws.Range("I135").Value = 1062.6818
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9564.136200000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7029.136200000001
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2329.65
$ws.Range("I137").Value = 1678.7
$ws.Range("J137").Value = 2492.3875
$ws.Range("K137").Value = 5036.1
$ws.Range("L137").Value = 7477.162499999999
$ws.Range("M137").Value = -2486.1
$ws.Range("N137").Value = -12577.1625
$ws.Range("H138").Value = 3367.491
$ws.Range("J138").Value = 4031.3635
$ws.Range("L138").Value = 12094.0905
$ws.Range("N138").Value = -22374.0905
$ws.Range("H141").Value = 3718.5
$ws.Range("I141").Value = 3718.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11155.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5975.5
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 103904.56
$ws.Range("I32").Value = 118661.914
$ws.Range("J32").Value = 7414.154
$ws.Range("K32").Value = 118661.914
$ws.Range("L32").Value = 7414.154
$ws.Range("M32").Value = -118374.914
$ws.Range("N32").Value = -7988.154
$ws.Range("H45").Value = 3210.1538
$ws.Range("I45").Value = 3064.75
$ws.Range("K45").Value = 3064.75
$ws.Range("M45").Value = -2687.75
$ws.Range("H61").Value = 1251433.5
$ws.Range("I61").Value = 715709.9
$ws.Range("K61").Value = 715709.9
$ws.Range("M61").Value = -715497.9
$ws.Range("H74").Value = 2764.4736
$ws.Range("I74").Value = 2925.5
$ws.Range("J74").Value = 2313.6
$ws.Range("K74").Value = 2925.5
$ws.Range("L74").Value = 2313.6
$ws.Range("M74").Value = -2051.5
$ws.Range("N74").Value = -4061.6
$ws.Range("H77").Value = 2764.4736
$ws.Range("I77").Value = 2925.5
$ws.Range("J77").Value = 2313.6
$ws.Range("K77").Value = 14627.5
$ws.Range("L77").Value = 11568
$ws.Range("M77").Value = -10259.5
$ws.Range("N77").Value = -20304
$ws.Range("H136").Value = 1251433.5
$ws.Range("I136").Value = 715709.9
$ws.Range("K136").Value = 2147129.7
$ws.Range("M136").Value = -2144579.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3189.3635
$ws.Range("I94").Value = 1760.375
$ws.Range("K94").Value = 1760.375
$ws.Range("M94").Value = -1309.375
$ws.Range("H107").Value = 13170130
$ws.Range("I107").Value = 8190.9653
$ws.Range("J107").Value = 55580820
$ws.Range("K107").Value = 8190.9653
$ws.Range("L107").Value = 55580820
$ws.Range("M107").Value = -6270.9653
$ws.Range("N107").Value = -55584660
$ws.Range("H134").Value = 3755.2856
$ws.Range("I134").Value = 2796.6667
$ws.Range("K134").Value = 8390.000100000001
$ws.Range("M134").Value = -5855.000100000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1714.5333
$ws.Range("I22").Value = 583.3333
$ws.Range("J22").Value = 1997.3334
$ws.Range("K22").Value = 583.3333
$ws.Range("L22").Value = 1997.3334
$ws.Range("M22").Value = -233.3333
$ws.Range("N22").Value = -2697.3334
$ws.Range("H31").Value = 2663.89
$ws.Range("J31").Value = 2889.573
$ws.Range("L31").Value = 2889.573
$ws.Range("N31").Value = -3479.573
$ws.Range("H34").Value = 2663.89
$ws.Range("J34").Value = 2889.573
$ws.Range("L34").Value = 2889.573
$ws.Range("N34").Value = -3293.573
$ws.Range("H58").Value = 2430.3076
$ws.Range("I58").Value = 1998.3334
$ws.Range("J58").Value = 2559.9
$ws.Range("K58").Value = 1998.3334
$ws.Range("L58").Value = 2559.9
$ws.Range("M58").Value = -1795.3334
$ws.Range("N58").Value = -2965.9
$ws.Range("H94").Value = 3227
$ws.Range("J94").Value = 2759.2
$ws.Range("L94").Value = 2759.2
$ws.Range("N94").Value = -3661.2
$ws.Range("H105").Value = 1700.8334
$ws.Range("I105").Value = 951.25
$ws.Range("J105").Value = 3200
$ws.Range("K105").Value = 951.25
$ws.Range("L105").Value = 3200
$ws.Range("M105").Value = 795.75
$ws.Range("N105").Value = -6694
$ws.Range("H132").Value = 2109.4644
$ws.Range("I132").Value = 1964.0385
$ws.Range("K132").Value = 5892.1155
$ws.Range("M132").Value = -3362.1155
$ws.Range("H136").Value = 2430.3076
$ws.Range("I136").Value = 1998.3334
$ws.Range("J136").Value = 2559.9
$ws.Range("K136").Value = 5995.0002
$ws.Range("L136").Value = 7679.700000000001
$ws.Range("M136").Value = -3445.0002
$ws.Range("N136").Value = -12779.7

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 67.1579
$ws.Range("I12").Value = 44.5
$ws.Range("K12").Value = 133.5
$ws.Range("M12").Value = 39.5
$ws.Range("H14").Value = 683.2353000000001
$ws.Range("I14").Value = 683.2353000000001
$ws.Range("K14").Value = 2049.7059
$ws.Range("M14").Value = -1876.7059
$ws.Range("H68").Value = 1853.0625
$ws.Range("I68").Value = 1083.3334
$ws.Range("J68").Value = 2030.6923
$ws.Range("K68").Value = 3250.0002
$ws.Range("L68").Value = 6092.0769
$ws.Range("M68").Value = -2439.0002
$ws.Range("N68").Value = -7714.0769
$ws.Range("H71").Value = 1853.0625
$ws.Range("I71").Value = 1083.3334
$ws.Range("J71").Value = 2030.6923
$ws.Range("K71").Value = 9750.000599999999
$ws.Range("L71").Value = 18276.2307
$ws.Range("M71").Value = -5694.000599999999
$ws.Range("N71").Value = -26388.2307

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9999
$ws.Range("J5").Value = 9999
$ws.Range("L5").Value = 9999
$ws.Range("N5").Value = -10223
$ws.Range("H70").Value = 13465.75
$ws.Range("I70").Value = 12550.4
$ws.Range("K70").Value = 12550.4
$ws.Range("M70").Value = -12280.4
$ws.Range("H73").Value = 13465.75
$ws.Range("I73").Value = 12550.4
$ws.Range("K73").Value = 12550.4
$ws.Range("M73").Value = -11614.4
$ws.Range("H97").Value = 1307.5333
$ws.Range("I97").Value = 1346.4584
$ws.Range("J97").Value = 1151.8334
$ws.Range("K97").Value = 1346.4584
$ws.Range("L97").Value = 1151.8334
$ws.Range("M97").Value = -850.4584
$ws.Range("N97").Value = -2143.8334
$ws.Range("H122").Value = 3609.4
$ws.Range("I122").Value = 3238.111
$ws.Range("J122").Value = 4166.3335
$ws.Range("K122").Value = 9714.332999999999
$ws.Range("L122").Value = 12499.0005
$ws.Range("M122").Value = -7264.332999999999
$ws.Range("N122").Value = -17399.0005
$ws.Range("H132").Value = 629884.5
$ws.Range("I132").Value = 1004464.1
$ws.Range("K132").Value = 3013392.3
$ws.Range("M132").Value = -3010862.3

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62503436
$ws.Range("J7").Value = 3928.2856
$ws.Range("L7").Value = 3928.2856
$ws.Range("N7").Value = -4152.2856
$ws.Range("H55").Value = 514.2632
$ws.Range("I55").Value = 561.1539
$ws.Range("J55").Value = 412.66666
$ws.Range("K55").Value = 561.1539
$ws.Range("L55").Value = 412.66666
$ws.Range("M55").Value = -388.1539
$ws.Range("N55").Value = -758.66666
$ws.Range("H93").Value = 1019.439
$ws.Range("I93").Value = 1054.75
$ws.Range("J93").Value = 893.8889
$ws.Range("K93").Value = 1054.75
$ws.Range("L93").Value = 893.8889
$ws.Range("M93").Value = 193.25
$ws.Range("N93").Value = -3389.8889
$ws.Range("H122").Value = 3949.9546
$ws.Range("I122").Value = 3545.4546
$ws.Range("J122").Value = 4354.4546
$ws.Range("K122").Value = 10636.3638
$ws.Range("L122").Value = 13063.3638
$ws.Range("M122").Value = -8186.363799999999
$ws.Range("N122").Value = -17963.3638
$ws.Range("H126").Value = 62503436
$ws.Range("J126").Value = 3928.2856
$ws.Range("L126").Value = 11784.8568
$ws.Range("N126").Value = -16724.8568

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13339530
$ws.Range("I81").Value = 4498.5
$ws.Range("K81").Value = 8997
$ws.Range("M81").Value = -7936
$ws.Range("H84").Value = 13339530
$ws.Range("I84").Value = 4498.5
$ws.Range("K84").Value = 44985
$ws.Range("M84").Value = -39681
$ws.Range("H122").Value = 5120.727
$ws.Range("I122").Value = 3991
$ws.Range("K122").Value = 11973
$ws.Range("M122").Value = -9523
